$wb = $excel.ActiveWorkbook

# Update zh-cn sheet: refresh handoff/handback datetimes for the
# 0bb40168-... file (row 2) as a result of re-generating the handback report.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 20:51:41"
$wsZhCn.Range("H2").Value = "2016-03-22 20:52:05"

# Update de-de sheet: refresh handoff/handback datetimes for the
# 0bb40168-... file (row 2) as a result of re-generating the handback report.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 20:51:45"
$wsDeDe.Range("H2").Value = "2016-03-22 20:52:11"
